$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data (refreshed snapshot).
# Numeric-looking price strings are forced to Text format so values such as
# "218.20" or "1.004" keep their exact original formatting (trailing zeros,
# three-decimal-group layout, etc.) instead of being auto-converted to numbers.

$ws.Range('D2').Value = '26.252.92'
$ws.Range('E2').Value = '  +0.12%  '

$ws.Range('D3').Value = '1.686.22'
$ws.Range('E3').Value = '  +0.79%  '

$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '218.20'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.24%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.5238'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +2.89%  '

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.06%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2693'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +1.62%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06434'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +1.74%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '22.02'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +0.97%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07454'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +0.91%  '

$ws.Range('D12').Value = '1.686.60'
$ws.Range('E12').Value = '  +0.78%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.532'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.30%  '

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.5854'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +1.58%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.000008555'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +0.01%  '

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '64.66'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.44%  '

$ws.Range('D17').Value = '26.309.81'
$ws.Range('E17').Value = '  +0.02%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '4.969'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.44%  '

$ws.Range('E19').Value = '  -0.15%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '10.84'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.00%  '

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '190.72'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +1.68%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.235'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.31%  '

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.11%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '145.28'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +1.78%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.1247'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +6.37%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '7.650'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +0.52%  '

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '15.83'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.71%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.06724'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +16.18%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.336'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.85%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.321'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.70%  '

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '3.596'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +2.31%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.565'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +1.68%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.663'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -0.45%  '

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.025'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +1.88%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.6187'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +3.36%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.374'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.18%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '2.706'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +2.24%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '6.287'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +6.21%  '

$ws.Range('D39').Value = '1.100.66'
$ws.Range('E39').Value = '  -0.06%  '

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.01609'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.75%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.8735'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.34%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.014'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +0.97%  '

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '101.03'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +1.54%  '

$ws.Range('D44').Value = '1.834.57'
$ws.Range('E44').Value = '  +0.77%  '

$ws.Range('E45').Value = '  +0.73%  '

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '56.91'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +1.10%  '

$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.005'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.13%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '8.133'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.84%  '

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.05256'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +1.05%  '

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.4291'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.52%  '

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '6.015'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +3.14%  '

Write-Output "Updated cryptos list"